$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1676554327.812017
$ws.Range("F2").Value = 457242089.4032773
$ws.Range("G2").Value = 2286210447.016387
$ws.Range("E3").Value = 4874074.231378612
$ws.Range("F3").Value = 1827777.83676698
$ws.Range("G3").Value = 9138889.1838349
$ws.Range("E4").Value = 1748273275.282674
$ws.Range("F4").Value = 476801802.3498201
$ws.Range("G4").Value = 2384009011.7491
$ws.Range("E5").Value = 4962811.308426417
$ws.Range("F5").Value = 1861054.240659907
$ws.Range("G5").Value = 9305271.203299535
$ws.Range("E6").Value = 2389772841.46038
$ws.Range("F6").Value = 651756229.4891943
$ws.Range("G6").Value = 3258781147.445971
$ws.Range("E7").Value = 6697133.033314341
$ws.Range("F7").Value = 2511424.887492878
$ws.Range("G7").Value = 12557124.43746439
$ws.Range("E8").Value = 2554361971.929386
$ws.Range("F8").Value = 696644174.1625597
$ws.Range("G8").Value = 3483220870.812799
$ws.Range("E9").Value = 7061336.111747485
$ws.Range("F9").Value = 2648001.041905307
$ws.Range("G9").Value = 13240005.20952654
$ws.Range("E10").Value = 2768713698.614011
$ws.Range("F10").Value = 755103735.9856393
$ws.Range("G10").Value = 3775518679.928197
$ws.Range("E11").Value = 7287355.023735483
$ws.Range("F11").Value = 2732758.133900806
$ws.Range("G11").Value = 13663790.66950403
$ws.Range("E12").Value = 2922421891.706268
$ws.Range("F12").Value = 797024152.2835275
$ws.Range("G12").Value = 3985120761.417638
$ws.Range("E13").Value = 7356243.281599058
$ws.Range("F13").Value = 2758591.230599647
$ws.Range("G13").Value = 13792956.15299824
$ws.Range("E14").Value = 3081982928.470347
$ws.Range("F14").Value = 840540798.6737309
$ws.Range("G14").Value = 4202703993.368654
$ws.Range("E15").Value = 7444045.18904587
$ws.Range("F15").Value = 2791516.945892202
$ws.Range("G15").Value = 13957584.72946101
$ws.Range("E16").Value = 3357934556.155107
$ws.Range("F16").Value = 915800333.4968474
$ws.Range("G16").Value = 4579001667.484237
$ws.Range("E17").Value = 7640031.331508775
$ws.Range("F17").Value = 2865011.749315791
$ws.Range("G17").Value = 14325058.74657895
$ws.Range("E19").Value = 7809483.843159612
$ws.Range("F19").Value = 2928556.441184855
$ws.Range("G19").Value = 14642782.20592427
